$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Adresse2" header in D1 and give it the same (bold/bordered) style as
# the other header cells in row 1.
$ws.Range("D1").Value = "Adresse2"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Remove the old "Adresse" values from column B - only the header stays.
$ws.Range("B2:B4").ClearContents()

# Size the new column to fit its (future) contents.
$ws.Columns.Item(4).ColumnWidth = 51.140625

# Restore the active selection similar to the authored workbook.
$ws.Range("B4").Select()

$wb.Save()
